$wb = $excel.ActiveWorkbook

# Map of old text -> new text, applied to every worksheet in the workbook.
$replacements = @{
    "bleu" = "noir"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($replacements.ContainsKey($val)) {
                $cell.Value2 = $replacements[$val]
            }
        }
    }
}
